# Apply config-file update: refresh recipient lists (To/CC) on the Config sheet,
# adjust row heights for the wrapped addresses, and remove the now-unused
# standalone "lester.rollan@..." shared string by overwriting the cells so
# they point at the already-existing To/CC list strings (shared with
# ConfigOptions rows 11/12).

$wb = $excel.ActiveWorkbook

$configSheet = $wb.Worksheets.Item("Config")
$optionsSheet = $wb.Worksheets.Item("ConfigOptions")

# --- Config sheet: rows 10 (RecipientTo) and 11 (RecipientCC) ---
$configSheet.Activate()

$configSheet.Range("B10").Value = "sam.tecson@lexisnexisrisk.com; joavic.quisano@lexisnexisrisk.com; david.villasoto@lexisnexisrisk.com"
$configSheet.Range("B11").Value = "lester.rollan@lexisnexisrisk.com; dindee.galindo@lexisnexisrisk.com; jesriel.tolentino@lexisnexisrisk.com; jhoanna.talle@lexisnexisrisk.com; paul.fabro@lexisnexisrisk.com; judy.cotaoco@lexisnexisrisk.com"

$configSheet.Rows.Item(10).RowHeight = 30
$configSheet.Rows.Item(11).RowHeight = 45

$configSheet.Range("B10:B11").Select()

# --- ConfigOptions sheet: update selection only (B39:B40) ---
$optionsSheet.Activate()
$optionsSheet.Range("B39:B40").Select()

$configSheet.Activate()
